$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'Duel Decks Sorin vs. Tibalt Tokens (TDDK)',
    'Aethersnipe',
    'Air Elemental',
    'Ancestral Vision',
    'Bottle Gnomes',
    'Brine Elemental',
    'Chandra Nalaar',
    'Chartooth Cougar',
    'Condescend',
    'Cone of Flame',
    'Counterspell',
    'Daze',
    'Demonfire',
    'Errant Ephemeron',
    'Fact or Fiction',
    'Fathom Seer',
    'Fireball',
    'Fireblast',
    'Firebolt',
    'Fireslinger',
    'Flame Javelin',
    'Flamekin Brawler',
    'Flametongue Kavu',
    'Flamewave Invoker',
    'Fledgling Mawcor',
    'Furnace Whelp',
    'Guile',
    'Gush',
    'Hostility',
    'Incinerate',
    'Ingot Chewer',
    'Inner-Flame Acolyte',
    'Island',
    'Island',
    'Island',
    'Island',
    'Jace Beleren',
    'Keldon Megaliths',
    'Magma Jet',
    'Man-o''-War',
    'Martyr of Frost',
    'Mind Stone',
    'Mountain',
    'Mountain',
    'Mountain',
    'Mountain',
    'Mulldrifter',
    'Ophidian',
    'Oxidda Golem',
    'Pyre Charger',
    'Quicksilver Dragon',
    'Rakdos Pit Dragon',
    'Repulse',
    'Riftwing Cloudskate',
    'Seal of Fire',
    'Slith Firewalker',
    'Soulbright Flamekin',
    'Spire Golem',
    'Terrain Generator',
    'Voidmage Apprentice',
    'Wall of Deceit',
    'Waterspout Djinn',
    'Willbender'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

Write-Output ("Wrote {0} rows" -f $values.Length)
